{"js": "// Replace the date line and every two-digit \u00d7 two-digit multiplication\n// prompt in the worksheet with the new values from the target revision.\nconst replacements = [\n  [\"2024-05-18 Saturday\", \"2024-05-19 Sunday\"],\n  [\"30\u00d719=\", \"22\u00d782=\"],\n  [\"82\u00d749=\", \"27\u00d787=\"],\n  [\"49\u00d762=\", \"19\u00d742=\"],\n  [\"66\u00d789=\", \"88\u00d732=\"],\n  [\"32\u00d735=\", \"91\u00d765=\"],\n  [\"52\u00d754=\", \"28\u00d768=\"],\n  [\"89\u00d730=\", \"47\u00d751=\"],\n  [\"39\u00d772=\", \"70\u00d727=\"],\n  [\"44\u00d725=\", \"82\u00d771=\"],\n  [\"88\u00d730=\", \"49\u00d786=\"],\n  [\"76\u00d797=\", \"79\u00d782=\"],\n  [\"26\u00d770=\", \"70\u00d787=\"],\n  [\"72\u00d714=\", \"74\u00d782=\"],\n  [\"50\u00d769=\", \"45\u00d758=\"],\n  [\"13\u00d747=\", \"74\u00d788=\"],\n  [\"54\u00d788=\", \"55\u00d748=\"],\n  [\"11\u00d783=\", \"59\u00d778=\"],\n  [\"82\u00d712=\", \"58\u00d756=\"],\n  [\"55\u00d721=\", \"58\u00d758=\"],\n  [\"99\u00d723=\", \"77\u00d732=\"],\n  [\"42\u00d734=\", \"72\u00d736=\"],\n  [\"91\u00d726=\", \"33\u00d789=\"],\n  [\"38\u00d748=\", \"24\u00d794=\"],\n  [\"67\u00d784=\", \"53\u00d797=\"],\n  [\"32\u00d715=\", \"41\u00d769=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and every two-digit \u00d7 two-digit multiplication\n# prompt in the worksheet with the new values from the target revision.\n$d = $word.ActiveDocument\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n$replacements = @(\n    @(\"2024-05-18 Saturday\", \"2024-05-19 Sunday\"),\n    @(\"30\u00d719=\", \"22\u00d782=\"),\n    @(\"82\u00d749=\", \"27\u00d787=\"),\n    @(\"49\u00d762=\", \"19\u00d742=\"),\n    @(\"66\u00d789=\", \"88\u00d732=\"),\n    @(\"32\u00d735=\", \"91\u00d765=\"),\n    @(\"52\u00d754=\", \"28\u00d768=\"),\n    @(\"89\u00d730=\", \"47\u00d751=\"),\n    @(\"39\u00d772=\", \"70\u00d727=\"),\n    @(\"44\u00d725=\", \"82\u00d771=\"),\n    @(\"88\u00d730=\", \"49\u00d786=\"),\n    @(\"76\u00d797=\", \"79\u00d782=\"),\n    @(\"26\u00d770=\", \"70\u00d787=\"),\n    @(\"72\u00d714=\", \"74\u00d782=\"),\n    @(\"50\u00d769=\", \"45\u00d758=\"),\n    @(\"13\u00d747=\", \"74\u00d788=\"),\n    @(\"54\u00d788=\", \"55\u00d748=\"),\n    @(\"11\u00d783=\", \"59\u00d778=\"),\n    @(\"82\u00d712=\", \"58\u00d756=\"),\n    @(\"55\u00d721=\", \"58\u00d758=\"),\n    @(\"99\u00d723=\", \"77\u00d732=\"),\n    @(\"42\u00d734=\", \"72\u00d736=\"),\n    @(\"91\u00d726=\", \"33\u00d789=\"),\n    @(\"38\u00d748=\", \"24\u00d794=\"),\n    @(\"67\u00d784=\", \"53\u00d797=\"),\n    @(\"32\u00d715=\", \"41\u00d769=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
